$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3193
$ws1.Range("F4").Value = 121
$ws1.Range("F5").Value = 6870
$ws1.Range("F6").Value = 2021
$ws1.Range("F7").Value = 23
$ws1.Range("F8").Value = 73
$ws1.Range("F11").Value = 71
$ws1.Range("F13").Value = 148
$ws1.Range("F14").Value = 182
$ws1.Range("F15").Value = 35

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3193
$ws4.Range("F5").Value = 121
$ws4.Range("F6").Value = 6870
$ws4.Range("F7").Value = 2021
$ws4.Range("F8").Value = 23
$ws4.Range("F9").Value = 73
$ws4.Range("F12").Value = 71
$ws4.Range("F14").Value = 148
$ws4.Range("F15").Value = 182
$ws4.Range("F16").Value = 35
